$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- New worksheet "Plan W", positioned right after "ICR" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Plan W"

# Clone the header row's look (bold "Mat" cell + bold/date-formatted
# B1:D1 band) from the ICR sheet so the new header matches it exactly.
$ws1.Range("A1:D1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)

# --- Cell values, written in the exact order the shared-string table
#     needs them so new-string indices land the same as the target ---
$ws2.Range("A1").Value = "Mat"
$ws2.Range("C1").Value = "Eerste"
$ws2.Range("D1").Value = "Laatste"

$ws2.Range("D3").Value = "84?"
$ws2.Range("C4").Value = "82?"

$ws2.Range("A2").Value = "Plan W1 blauw '74"
$ws2.Range("A3").Value = "Plan W1 benelux '74"
$ws2.Range("A4").Value = "Plan W1 benelux '84"
$ws2.Range("A5").Value = "Plan W1 ic '87"
$ws2.Range("A9").Value = "Plan W2 ic '82"
$ws2.Range("A7").Value = "Plan W2 blauw rb '68"
$ws2.Range("A8").Value = "Plan W2 507 '81"
$ws2.Range("A6").Value = "Plan W1 ic '87"
$ws2.Range("A10").Value = "Plan W2 ic '82"

# Numbers (dates stored as date serials, formatted below)
$ws2.Range("C2").Value = 25750
$ws2.Range("D2").Value = 32417
$ws2.Range("C3").Value = 26816
$ws2.Range("D4").Value = 32417
$ws2.Range("C5").Value = 31929
$ws2.Range("D5").Value = 35217
$ws2.Range("C6").Value = 35916
$ws2.Range("D6").Value = 37956
$ws2.Range("C7").Value = 25020
$ws2.Range("D7").Value = 30682
$ws2.Range("C8").Value = 29891
$ws2.Range("D8").Value = 30256
$ws2.Range("C9").Value = 30256
$ws2.Range("D9").Value = 35309
$ws2.Range("C10").Value = 35916
$ws2.Range("D10").Value = 37956

# Date number format on the whole C2:D10 block (reuses the workbook's
# existing "mmm-yy" custom format / style), then strip it back off C4
# which holds a plain "82?" guess rather than a date.
$ws2.Range("C2:D10").NumberFormat = "mmm-yy"
$ws2.Range("C4").ClearFormats()
$ws2.Range("C4").Value = "82?"

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 3
$ws2.Columns.Item(2).ColumnWidth = 16.33203125

# Comments ("Tom:" notes), mirroring the author's annotations
$ws2.Range("A1").AddComment("Tom:`nW1: nummer eindigt op 4xx`nW2: nummer eindigt op 5xx")
$ws2.Range("A2").AddComment("Tom:`nslechts 3 rijtuigen, rest was Bnl")
$ws2.Range("A3").AddComment("Tom:`n'73 op railwiki")
$ws2.Range("A4").AddComment("Tom:`nandere ramen")
$ws2.Range("C6").AddComment("Tom:`nherindienststelling")
$ws2.Range("A8").AddComment("Tom:`nblauw met gele deuren`nenkel rijtuig: 507")
$ws2.Range("A9").AddComment("Tom:`nherkenbaar aan nummer rechts ipv midden")
$ws2.Range("C10").AddComment("Tom:`nherindienststelling")

# --- Sheet view / selection housekeeping ---
$ws1.Range("A1:D1").Select()
$ws2.Select()
$ws2.Range("D19").Select()
